$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column header
$ws.Range("D1").Value = "FilePath"

# Add new rows of data
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "chill"
$ws.Range("C6").Value = 535.6

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "water"
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = "D:\work-place\flutter apps\projects\sharp\BisleriumCafeBackend\fyp-document\fyp\coffee\coffee-image\2024-01-08\1704735419442-5ce111d5-dd7c-41f3-b432-7abff9a14dd6.jpg"

# Update selection to match final state
$ws.Range("D7").Select()
